$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Row 2 (OTROS)
$ws.Range("D2").Value = 6005.41
$ws.Range("E2").Value = -6005.41

# Row 3 (PORCELANATO)
$ws.Range("D3").Value = 14195.66
$ws.Range("E3").Value = -472.3199999999997
$ws.Range("F3").Value = 1.034417277426632

# Row 4 (TOTAL)
$ws.Range("D4").Value = 20201.07
$ws.Range("E4").Value = -6477.73
$ws.Range("F4").Value = 1.472022845750379
